$wb = $excel.ActiveWorkbook

# Sheets: Germany, Belgium, Czech
$wsGermany = $wb.Worksheets.Item("Germany")
$wsBelgium = $wb.Worksheets.Item("Belgium")
$wsCzech   = $wb.Worksheets.Item("Czech")

# "Other Nodes" list: replace "MOXA Node" entry with "BACnet Interface" on every market sheet (cell A10)
$wsGermany.Range("A10").Value = "BACnet Interface"
$wsBelgium.Range("A10").Value = "BACnet Interface"
$wsCzech.Range("A10").Value = "BACnet Interface"

# Update each sheet's remembered selection; select Germany's range last so it ends up
# the active/selected sheet (tabSelected) when the workbook is saved.
$wsBelgium.Range("A10").Select() | Out-Null
$wsCzech.Range("A10").Select() | Out-Null
$wsGermany.Range("A8").Select() | Out-Null

$wb.Save()
